$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.725.97"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "'1.919.43"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'240.19"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.4928"
$ws.Range("D8").Value = "'0.2986"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").Value = "'0.06781"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "'1.944.73"
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D11").Value = "'17.28"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "'0.07355"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "'5.188"
$ws.Range("E13").Value = "  +2.66%  "
$ws.Range("D14").Value = "'89.04"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "'0.6746"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "'30.710.68"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'0.000007955"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "'13.58"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "'2.160.83"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'5.356"
$ws.Range("E21").Value = "  +11.15%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'201.92"
$ws.Range("E23").Value = "  +5.20%  "
$ws.Range("D24").Value = "'6.316"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("D25").Value = "'9.678"
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("D26").Value = "'166.45"
$ws.Range("E26").Value = "  +6.77%  "
$ws.Range("D27").Value = "'18.92"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("E28").Value = "  +3.44%  "
$ws.Range("D29").Value = "'1.477"
$ws.Range("E29").Value = "  +5.47%  "
$ws.Range("D30").Value = "'4.376"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").Value = "'0.09172"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "'4.075"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("D33").Value = "'0.05312"
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("D34").Value = "'0.7440"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "'1.122"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").Value = "'2.729"
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").Value = "'0.01841"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").Value = "'2.725"
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("D39").Value = "'0.9246"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("D40").Value = "'2.086"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").Value = "'75.57"
$ws.Range("E41").Value = "  +30.98%  "
$ws.Range("D42").Value = "'0.4476"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").Value = "'5.987"
$ws.Range("E43").Value = "  +4.21%  "
$ws.Range("D44").Value = "'107.13"
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "'0.1393"
$ws.Range("E46").Value = "  +3.72%  "
$ws.Range("D47").Value = "'7.666"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("D48").Value = "'35.85"
$ws.Range("E48").Value = "  +6.55%  "
$ws.Range("D49").Value = "'9.108"
$ws.Range("E49").Value = "  +4.01%  "
$ws.Range("D50").Value = "'0.05887"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("E51").Value = "  +2.39%  "
